$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '30.591.25'

$ws.Range("D3").Value = '1.877.14'
$ws.Range("E3").Value = '  -0.18%  '

Set-TextCell "D4" '0.9997'
$ws.Range("E4").Value = '  -0.01%  '

Set-TextCell "D5" '238.89'
$ws.Range("E5").Value = '  +0.61%  '

Set-TextCell "D6" '0.9996'
$ws.Range("E6").Value = '  +0.01%  '

Set-TextCell "D7" '0.4801'
$ws.Range("E7").Value = '  -0.48%  '

Set-TextCell "D8" '0.2835'
$ws.Range("E8").Value = '  -1.95%  '

Set-TextCell "D9" '0.06527'
$ws.Range("E9").Value = '  -0.90%  '

$ws.Range("D10").Value = '1.962.57'
$ws.Range("E10").Value = '  +4.45%  '

Set-TextCell "D11" '0.07463'
$ws.Range("E11").Value = '  +0.98%  '

Set-TextCell "D12" '16.62'
$ws.Range("E12").Value = '  -1.81%  '

Set-TextCell "D13" '5.100'
$ws.Range("E13").Value = '  -1.75%  '

Set-TextCell "D14" '88.82'
$ws.Range("E14").Value = '  +0.90%  '

Set-TextCell "D15" '0.6651'
$ws.Range("E15").Value = '  +0.79%  '

$ws.Range("D16").Value = '30.567.19'
$ws.Range("E16").Value = '  +0.92%  '

$ws.Range("E17").Value = '  -2.11%  '

Set-TextCell "D18" '0.9996'
$ws.Range("E18").Value = '  -0.02%  '

Set-TextCell "D19" '0.000007612'
$ws.Range("E19").Value = '  -1.51%  '

# Row 20: coin name/link swap
$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.163.55'
$ws.Range("E20").Value = '  +1.21%  '

# Row 21: coin name/link swap
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell "D21" '230.44'
$ws.Range("E21").Value = '  +17.60%  '

Set-TextCell "D22" '5.307'
$ws.Range("E22").Value = '  -3.04%  '

Set-TextCell "D23" '0.9993'
$ws.Range("E23").Value = '  -0.07%  '

Set-TextCell "D24" '6.210'
$ws.Range("E24").Value = '  +0.84%  '

Set-TextCell "D25" '9.318'
$ws.Range("E25").Value = '  -1.27%  '

Set-TextCell "D26" '167.30'
$ws.Range("E26").Value = '  +2.40%  '

Set-TextCell "D28" '1.953'
$ws.Range("E28").Value = '  +1.33%  '

Set-TextCell "D29" '1.455'
$ws.Range("E29").Value = '  +0.98%  '

Set-TextCell "D30" '0.09544'
$ws.Range("E30").Value = '  +4.34%  '

Set-TextCell "D31" '4.314'

Set-TextCell "D32" '4.033'
$ws.Range("E32").Value = '  -0.49%  '

Set-TextCell "D33" '0.05031'
$ws.Range("E33").Value = '  -0.31%  '

Set-TextCell "D34" '1.214'
$ws.Range("E34").Value = '  +6.42%  '

Set-TextCell "D35" '0.7489'
$ws.Range("E35").Value = '  +1.00%  '

Set-TextCell "D36" '2.712'
$ws.Range("E36").Value = '  +0.21%  '

Set-TextCell "D37" '0.01848'
$ws.Range("E37").Value = '  +0.33%  '

Set-TextCell "D38" '2.623'
$ws.Range("E38").Value = '  -0.27%  '

Set-TextCell "D39" '2.077'
$ws.Range("E39").Value = '  +0.27%  '

Set-TextCell "D40" '0.9104'
$ws.Range("E40").Value = '  -0.63%  '

Set-TextCell "D41" '105.94'
$ws.Range("E41").Value = '  -0.41%  '

# Row 42: coin name/link swap
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell "D42" '0.4279'
$ws.Range("E42").Value = '  -0.97%  '

# Row 43: coin name/link swap
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell "D43" '5.805'
$ws.Range("E43").Value = '  -1.08%  '

$ws.Range("E44").Value = '  +0.58%  '

Set-TextCell "D45" '7.486'
$ws.Range("E45").Value = '  -1.97%  '

Set-TextCell "D46" '64.57'
$ws.Range("E46").Value = '  -0.85%  '

Set-TextCell "D47" '0.1285'
$ws.Range("E47").Value = '  -4.74%  '

Set-TextCell "D48" '1.483'
$ws.Range("E48").Value = '  -5.53%  '

Set-TextCell "D49" '8.957'
$ws.Range("E49").Value = '  +0.74%  '

Set-TextCell "D50" '33.85'
$ws.Range("E50").Value = '  -0.95%  '

Set-TextCell "D51" '0.3884'
$ws.Range("E51").Value = '  +0.27%  '
